# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1151
$ws1.Range("F4").Value = 2607

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1151
$ws4.Range("F6").Value = 2607

$wb.Save()
